# [Fonds de solidarite] Add 2022-05-12 data
# Updates counts/amounts (nombre_aides, nombre_entreprises, montant_total)
# for a set of existing rows with newer cumulative figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 345532
$ws.Range("E10").Value = 1817701330

$ws.Range("C11").Value = 1265
$ws.Range("E11").Value = 46515184

$ws.Range("C13").Value = 187832
$ws.Range("E13").Value = 1165163461

$ws.Range("C88").Value = 71261
$ws.Range("E88").Value = 110287231

$ws.Range("C91").Value = 18845
$ws.Range("E91").Value = 75099352

$ws.Range("C93").Value = 16917
$ws.Range("E93").Value = 50434566

$ws.Range("C98").Value = 6285
$ws.Range("E98").Value = 19290013

$ws.Range("C100").Value = 9333
$ws.Range("E100").Value = 23704653

$ws.Range("C112").Value = 145226
$ws.Range("E112").Value = 716158923

$ws.Range("C119").Value = 8982
$ws.Range("E119").Value = 37087784

$ws.Range("C121").Value = 1306111
$ws.Range("E121").Value = 2274508151

$ws.Range("C129").Value = 633318
$ws.Range("E129").Value = 3426201653

$ws.Range("C130").Value = 4239
$ws.Range("E130").Value = 140350460

$ws.Range("C132").Value = 585588
$ws.Range("D132").Value = 90777
$ws.Range("E132").Value = 3459563006

$ws.Range("C139").Value = 76637
$ws.Range("E139").Value = 114129393

$ws.Range("C144").Value = 25065
$ws.Range("E144").Value = 92331999

$ws.Range("C145").Value = 72
$ws.Range("E145").Value = 6534296

$ws.Range("C146").Value = 7438
$ws.Range("E146").Value = 37681278

$ws.Range("C150").Value = 895
$ws.Range("E150").Value = 2015502

$ws.Range("C151").Value = 39920
$ws.Range("E151").Value = 60358891

$ws.Range("C154").Value = 18430
$ws.Range("E154").Value = 72554723

$ws.Range("C156").Value = 12394
$ws.Range("E156").Value = 40027556

$ws.Range("C186").Value = 236818
$ws.Range("E186").Value = 1189702070

$ws.Range("C194").Value = 18378
$ws.Range("E194").Value = 71345179

$ws.Range("C215").Value = 230252
$ws.Range("E215").Value = 408700746

$ws.Range("C221").Value = 135494
$ws.Range("E221").Value = 681816845

$ws.Range("C240").Value = 205897
$ws.Range("E240").Value = 1068613038
